# Actualización y re-estructuración del código
# - Rename row labels in column A (strip the ".mean" suffix from the
#   variable names used in the grouped-statistics table).
# - Refresh the computed statistics (mean/sd/min/25%/50%/75%/max) in
#   columns B:H for each variable row, rows 2-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: row labels, ".mean" suffix removed -------------------------
$labels = @{
    2  = "goldEarned"
    3  = "kills"
    4  = "deaths"
    5  = "assists"
    6  = "champExperience"
    7  = "player.WR"
    8  = "turretKills"
    9  = "totalMinionsKilled"
    10 = "totalTimeCCDealt"
    11 = "baronKills"
    12 = "dragonKills"
    13 = "totalDamageDealt"
    14 = "totalDamageTaken"
    15 = "totalDamageDealtToChampions"
    16 = "damageDealtToObjectives"
    17 = "goldEarnedPerMinute"
    18 = "visionScore"
}

foreach ($r in $labels.Keys) {
    $ws.Range("A$r").Value = $labels[$r]
}

# --- Columns B:H: updated statistics ---------------------------------------
$data = @{
    2  = @(10382.1815339248, 3852.53681543654, 500, 7544, 9997, 12896, 34146)
    3  = @(5.28342273177924, 4.32393789092404, 0, 2, 4, 8, 34)
    4  = @(4.68032641024989, 2.88871937527078, 0, 2, 4, 7, 21)
    5  = @(7.78218979586649, 5.86532339793899, 0, 3, 7, 11, 44)
    6  = @(11376.2405623205, 4295.2504418734, 0, 8228, 11205, 14214, 33483)
    7  = @(0.538627755039784, 0.0651390355652344, 0.285714285714286, 0.494623655913978, 0.536842105263158, 0.578313253012048, 0.863636363636364)
    8  = @(1.07204413147927, 1.38759051478975, 0, 0, 1, 2, 11)
    9  = @(119.788545135623, 89.7849912242387, 0, 28, 132, 195, 492)
    10 = @(243.634745913517, 353.129174670099, 0, 80, 161, 288, 6151)
    11 = @(0.102270127361009, 0.335723354528011, 0, 0, 0, 0, 3)
    12 = @(0.316445077153825, 0.737393054875219, 0, 0, 0, 0, 6)
    13 = @(121562.940361492, 81143.5594554608, 0, 57244, 115956, 172198, 639246)
    14 = @(20014.8134834888, 10341.2142431826, 0, 12494, 18396, 25772, 136330)
    15 = @(16520.6634466279, 11066.9406952063, 0, 8140, 14122, 22353, 117333)
    16 = @(10770.8425909449, 12394.2069318237, 0, 1894, 5949, 15392, 99712)
    17 = @(418.080040864473, 103.232457931167, 177.410071942446, 341.941747572816, 415.930949445129, 486.714801444043, 2500)
    18 = @(30.5680911101507, 23.1021717229153, 0, 16, 24, 37, 280)
}

$columns = @("B", "C", "D", "E", "F", "G", "H")

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $r).Value = $values[$i]
    }
}
